# The deck ships with two themes:
#   ppt/theme/theme1.xml  -> bound to the slide master (currently the
#                             "Integral" color scheme)
#   ppt/theme/theme2.xml  -> bound to the notes master (currently the
#                             default "Office Theme" color scheme)
#
# The authored edit swaps the two themes' contents (the slide master
# picks up the plain "Office Theme" colors, the notes master picks up
# the "Integral" colors). The font scheme and format scheme (fills,
# lines, effects) are byte-identical between the two themes already,
# so the only real content that moves is the 12-slot color scheme
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
#
# Recolor the slide master's theme (theme1.xml) from "Integral" to the
# stock "Office Theme" palette via the live ThemeColorScheme, which is
# the PowerPoint object-model surface for editing a:clrScheme in place
# without disturbing the rest of the theme (fonts/format scheme/names).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeColors[$i - 1]
}
